# Scheduled Kujata Profits runner: refresh Universalis market-price snapshots
# and recompute Leve crafting-profit columns (currentAveragePrice*, Leve*Price*,
# LeveProfit*) across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 3541.7334
$ws.Range("I6").Value = 5182.6
$ws.Range("J6").Value = 260
$ws.Range("K6").Value = 15547.8
$ws.Range("L6").Value = 780
$ws.Range("M6").Value = -15435.8
$ws.Range("N6").Value = -1004

# Row 9
$ws.Range("H9").Value = 145.16667
$ws.Range("I9").Value = 145.16667
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 145.16667
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = 23.83332999999999
$ws.Range("N9").ClearContents()

# Row 12
$ws.Range("H12").Value = 323.16666
$ws.Range("I12").Value = 287.8
$ws.Range("K12").Value = 287.8
$ws.Range("M12").Value = -117.8

# Row 19
$ws.Range("H19").Value = 3998.5
$ws.Range("I19").Value = 10000
$ws.Range("J19").Value = 2798.2
$ws.Range("K19").Value = 10000
$ws.Range("L19").Value = 2798.2
$ws.Range("M19").Value = -9825
$ws.Range("N19").Value = -3148.2

# Row 28
$ws.Range("H28").Value = 6519.9375
$ws.Range("I28").Value = 6920.8667
$ws.Range("K28").Value = 6920.8667
$ws.Range("M28").Value = -6435.8667

# Row 40
$ws.Range("H40").Value = 1759.1538
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 2061.5
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 2061.5
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2411.5

# Row 53
$ws.Range("H53").Value = 3502
$ws.Range("I53").Value = 4666
$ws.Range("J53").Value = 10
$ws.Range("K53").Value = 4666
$ws.Range("L53").Value = 10
$ws.Range("M53").Value = -4029
$ws.Range("N53").Value = -1284

# Row 92
$ws.Range("H92").Value = 1187.0834
$ws.Range("I92").Value = 1187.0834
$ws.Range("K92").Value = 1187.0834
$ws.Range("M92").Value = 60.91660000000002

# Row 98
$ws.Range("H98").Value = 1649.2142
$ws.Range("I98").Value = 1649.2142
$ws.Range("K98").Value = 1649.2142
$ws.Range("M98").Value = -151.2141999999999

# Row 122
$ws.Range("H122").Value = 1649.2142
$ws.Range("I122").Value = 1649.2142
$ws.Range("K122").Value = 4947.642599999999
$ws.Range("M122").Value = -2497.642599999999

# Row 125
$ws.Range("H125").Value = 4933.3335
$ws.Range("I125").Value = 4900
$ws.Range("K125").Value = 44100
$ws.Range("M125").Value = -41640

# Row 132
$ws.Range("H132").Value = 13339648
$ws.Range("I132").Value = 17547594
$ws.Range("K132").Value = 52642782
$ws.Range("M132").Value = -52640252

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1359.3
$ws.Range("I61").Value = 1288.1111
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1288.1111
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1076.1111
$ws.Range("N61").Value = -2424

# Row 74
$ws.Range("H74").Value = 741.91895
$ws.Range("I74").Value = 680.9355
$ws.Range("J74").Value = 1057
$ws.Range("K74").Value = 680.9355
$ws.Range("L74").Value = 1057
$ws.Range("M74").Value = 193.0645
$ws.Range("N74").Value = -2805

# Row 77
$ws.Range("H77").Value = 741.91895
$ws.Range("I77").Value = 680.9355
$ws.Range("J77").Value = 1057
$ws.Range("K77").Value = 3404.6775
$ws.Range("L77").Value = 5285
$ws.Range("M77").Value = 963.3224999999998
$ws.Range("N77").Value = -14021

# Row 132
$ws.Range("H132").Value = 3655.7896
$ws.Range("I132").Value = 3702.8333
$ws.Range("J132").Value = 3575.1428
$ws.Range("K132").Value = 11108.4999
$ws.Range("L132").Value = 10725.4284
$ws.Range("M132").Value = -8578.499899999999
$ws.Range("N132").Value = -15785.4284

# Row 136
$ws.Range("H136").Value = 1359.3
$ws.Range("I136").Value = 1288.1111
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3864.3333
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1314.3333
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1564.6923
$ws.Range("I107").Value = 1192.8
$ws.Range("J107").Value = 2804.3333
$ws.Range("K107").Value = 1192.8
$ws.Range("L107").Value = 2804.3333
$ws.Range("M107").Value = 727.2
$ws.Range("N107").Value = -6644.3333

# Row 134
$ws.Range("H134").Value = 10576.044
$ws.Range("I134").Value = 7336.1113
$ws.Range("J134").Value = 22239.8
$ws.Range("K134").Value = 22008.3339
$ws.Range("L134").Value = 66719.39999999999
$ws.Range("M134").Value = -19473.3339
$ws.Range("N134").Value = -71789.39999999999

# Row 135
$ws.Range("H135").Value = 38937.6
$ws.Range("J135").Value = 38937.6
$ws.Range("L135").Value = 38937.6
$ws.Range("N135").Value = -49077.6

# Row 140
$ws.Range("H140").Value = 22063.285
$ws.Range("J140").Value = 22127.45
$ws.Range("L140").Value = 22127.45
$ws.Range("N140").Value = -32487.45

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# Row 31
$ws.Range("H31").Value = 860.75806
$ws.Range("J31").Value = 1036.4546
$ws.Range("L31").Value = 1036.4546
$ws.Range("N31").Value = -1626.4546

# Row 34
$ws.Range("H34").Value = 860.75806
$ws.Range("J34").Value = 1036.4546
$ws.Range("L34").Value = 1036.4546
$ws.Range("N34").Value = -1440.4546

# Row 58
$ws.Range("H58").Value = 849.26086
$ws.Range("I58").Value = 863.6875
$ws.Range("J58").Value = 816.2857
$ws.Range("K58").Value = 863.6875
$ws.Range("L58").Value = 816.2857
$ws.Range("M58").Value = -660.6875
$ws.Range("N58").Value = -1222.2857

# Row 107
$ws.Range("H107").Value = 691.86365
$ws.Range("I107").Value = 569.75
$ws.Range("J107").Value = 761.6429000000001
$ws.Range("K107").Value = 569.75
$ws.Range("L107").Value = 761.6429000000001
$ws.Range("M107").Value = 1350.25
$ws.Range("N107").Value = -4601.6429

# Row 134
$ws.Range("H134").Value = 12822007
$ws.Range("I134").Value = 14494078
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 43482234
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -43479699
$ws.Range("N134").Value = -13470

# Row 136
$ws.Range("H136").Value = 849.26086
$ws.Range("I136").Value = 863.6875
$ws.Range("J136").Value = 816.2857
$ws.Range("K136").Value = 2591.0625
$ws.Range("L136").Value = 2448.8571
$ws.Range("M136").Value = -41.0625
$ws.Range("N136").Value = -7548.8571

# Row 141
$ws.Range("H141").Value = 33651
$ws.Range("J141").Value = 33601.145
$ws.Range("L141").Value = 33601.145
$ws.Range("N141").Value = -43961.145

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1365.3226
$ws.Range("I68").Value = 684.6429000000001
$ws.Range("K68").Value = 2053.9287
$ws.Range("M68").Value = -1242.9287

# Row 71
$ws.Range("H71").Value = 1365.3226
$ws.Range("I71").Value = 684.6429000000001
$ws.Range("K71").Value = 6161.7861
$ws.Range("M71").Value = -2105.7861

# Row 107
$ws.Range("H107").Value = 5104.5654
$ws.Range("J107").Value = 9996.182000000001
$ws.Range("L107").Value = 29988.546
$ws.Range("N107").Value = -33828.546

# Row 122
$ws.Range("H122").Value = 751.41174
$ws.Range("J122").Value = 811.06665
$ws.Range("L122").Value = 7299.59985
$ws.Range("N122").Value = -12199.59985

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 209.58824
$ws.Range("I2").Value = 185.625
$ws.Range("K2").Value = 185.625
$ws.Range("M2").Value = -72.625

# Row 45
$ws.Range("H45").Value = 38163
$ws.Range("J45").Value = 38163
$ws.Range("L45").Value = 38163
$ws.Range("N45").Value = -39281

# Row 62
$ws.Range("H62").Value = 16723.334

# Row 65
$ws.Range("H65").Value = 16723.334

# Row 70
$ws.Range("H70").Value = 23689290
$ws.Range("I70").Value = 35718600
$ws.Range("J70").Value = 16672192
$ws.Range("K70").Value = 35718600
$ws.Range("L70").Value = 16672192
$ws.Range("M70").Value = -35718330
$ws.Range("N70").Value = -16672732

# Row 73
$ws.Range("H73").Value = 23689290
$ws.Range("I73").Value = 35718600
$ws.Range("J73").Value = 16672192
$ws.Range("K73").Value = 35718600
$ws.Range("L73").Value = 16672192
$ws.Range("M73").Value = -35717664
$ws.Range("N73").Value = -16674064

# Row 126
$ws.Range("H126").Value = 2411.2727
$ws.Range("I126").Value = 1803.4286
$ws.Range("K126").Value = 5410.2858
$ws.Range("M126").Value = -2940.2858

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2550.6
$ws.Range("I7").Value = 2248
$ws.Range("K7").Value = 2248
$ws.Range("M7").Value = -2136

# Row 93
$ws.Range("H93").Value = 755.1111
$ws.Range("I93").Value = 698.8570999999999
$ws.Range("K93").Value = 698.8570999999999
$ws.Range("M93").Value = 549.1429000000001

# Row 106
$ws.Range("H106").Value = 27142
$ws.Range("J106").Value = 27142
$ws.Range("L106").Value = 27142
$ws.Range("N106").Value = -29666

# Row 126
$ws.Range("H126").Value = 2550.6
$ws.Range("I126").Value = 2248
$ws.Range("K126").Value = 6744
$ws.Range("M126").Value = -4274

# Row 136
$ws.Range("H136").Value = 2721.4443
$ws.Range("I136").Value = 2883.2856
$ws.Range("J136").Value = 2155
$ws.Range("K136").Value = 8649.856800000001
$ws.Range("L136").Value = 6465
$ws.Range("M136").Value = -6099.856800000001
$ws.Range("N136").Value = -11565

$ws = $wb.Worksheets.Item("WVR")
# Row 92
$ws.Range("H92").Value = 14610
$ws.Range("J92").Value = 14610
$ws.Range("L92").Value = 14610
$ws.Range("N92").Value = -19602

# Row 126
$ws.Range("H126").Value = 58825304
$ws.Range("I126").Value = 200000860
$ws.Range("J126").Value = 2158.3333
$ws.Range("K126").Value = 600002580
$ws.Range("L126").Value = 6474.999899999999
$ws.Range("M126").Value = -600000110
$ws.Range("N126").Value = -11414.9999

